$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 355, shifting existing rows 355:414 down to 356:415
$ws.Rows.Item(355).Insert()

# Populate the newly inserted row 355 with the new record's data
$ws.Cells.Item(355, 1).Value = 9
$ws.Cells.Item(355, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(355, 3).Value = "Metropolitana"
$ws.Cells.Item(355, 4).Value = 44776
$ws.Cells.Item(355, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(355, 5).Value = 13
$ws.Cells.Item(355, 6).Value = 100112044
$ws.Cells.Item(355, 7).Value = "Perejil"
$ws.Cells.Item(355, 8).Value = "Sin especificar"
$ws.Cells.Item(355, 9).Value = "Primera"
$ws.Cells.Item(355, 10).Value = 70
$ws.Cells.Item(355, 11).Value = 19000
$ws.Cells.Item(355, 12).Value = 20000
$ws.Cells.Item(355, 13).Value = 19500
$ws.Cells.Item(355, 14).Value = "`$/docena de atados"
$ws.Cells.Item(355, 15).Value = "Región Metropolitana"
$ws.Cells.Item(355, 16).Value = 6500
$ws.Cells.Item(355, 17).Value = 3
$ws.Cells.Item(355, 18).Value = "Hortaliza"
